# Add parameters to tests
$wb = $excel.ActiveWorkbook

# --- AddCustomerTest: append two more test-parameter rows ---
$wsAdd = $wb.Worksheets.Item("AddCustomerTest")

$wsAdd.Range("A3").Value = "Maria"
$wsAdd.Range("B3").Value = "Santos"
$wsAdd.Range("C3").Value = 654321
$wsAdd.Range("D3").Value = "Customer added successfully"

$wsAdd.Range("A4").Value = "Enzo"
$wsAdd.Range("B4").Value = "Alvez"
$wsAdd.Range("C4").Value = 123654
$wsAdd.Range("D4").Value = "Customer added successfully"

$wsAdd.Range("D5").Select() | Out-Null

# --- OpenAccountTest: swap the customer/currency sample values ---
$wsOpen = $wb.Worksheets.Item("OpenAccountTest")

$wsOpen.Range("A2").Value = "Joao Silva"
$wsOpen.Range("B2").Value = "Dollar"

$wsOpen.Range("A3").Select() | Out-Null

# --- Make AddCustomerTest the active sheet/tab ---
$wsAdd.Activate() | Out-Null
$wsAdd.Select() | Out-Null
$wsAdd.Range("D5").Select() | Out-Null
